$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contribuyente")

# Fix a handful of CCC (column H) data-entry errors before computing IBANs
$ws.Range("H77").Value = "20960056133231500000"
$ws.Range("H80").Value = "23658965274585223202"
$ws.Range("H83").Value = "20012541150023365233"
$ws.Range("H102").Value = "65645150865168448896"
$ws.Range("H103").Value = "26551681807651415636"
$ws.Range("H106").Value = "51556584221251000254"
$ws.Range("H128").Value = "62541122421110105611"

# Populate the new IBAN column (I) for every contributor row: country code + check digits + CCC
$ws.Range("I2").Value = "DK7331645124473461205164"
$ws.Range("I3").Value = "ES8265614874165615445616"
$ws.Range("I4").Value = "RO8832569523016220165156"
$ws.Range("I5").Value = "DE7424561937521546497521"
$ws.Range("I6").Value = "MC6436520125638451012515"
$ws.Range("I7").Value = "ES0721584976902154655487"
$ws.Range("I8").Value = "GR9420125003305201112544"
$ws.Range("I9").Value = "ES2821651484690980008984"
$ws.Range("I10").Value = "FI5620960043043554600000"
$ws.Range("I11").Value = "ES7921564975243245467995"
$ws.Range("I12").Value = "LT8032566221522587754554"
$ws.Range("I13").Value = "EE2023215465315456411515"
$ws.Range("I14").Value = "BE9400750184310702510000"
$ws.Range("I16").Value = "SM2125894363475485700145"
$ws.Range("I17").Value = "ES9596431245118150005156"
$ws.Range("I18").Value = "AT6825030000114574745458"
$ws.Range("I19").Value = "IT8915953684811254695203"
$ws.Range("I20").Value = "ES9020960043023096200000"
$ws.Range("I21").Value = "DK5800750184310702510000"
$ws.Range("I22").Value = "ES5023455254943263234457"
$ws.Range("I23").Value = "GR4920910936583000000000"
$ws.Range("I24").Value = "ES3720960043032159000000"
$ws.Range("I25").Value = "DE5512669681115112121210"
$ws.Range("I27").Value = "ES2956187775315550000651"
$ws.Range("I28").Value = "ES0425516848021156151054"
$ws.Range("I29").Value = "PT5764578946740051516490"
$ws.Range("I30").Value = "ES4534698752714600549403"
$ws.Range("I31").Value = "ES2766649444162310000255"
$ws.Range("I32").Value = "FR5623185484465641685100"
$ws.Range("I36").Value = "DE5021508149175421346497"
$ws.Range("I37").Value = "DE6721346154503164978451"
$ws.Range("I38").Value = "ES7225187786311225455548"
$ws.Range("I39").Value = "ES4723164897642213030615"
$ws.Range("I40").Value = "ES2396536214865214585214"
$ws.Range("I41").Value = "ES6885461325251978750005"
$ws.Range("I42").Value = "FI5024587946032003165464"
$ws.Range("I43").Value = "ES5020960043073071400000"
$ws.Range("I44").Value = "ES8220960043042158800000"
$ws.Range("I45").Value = "ES7521654587985156484454"
$ws.Range("I46").Value = "ES3251651681961210656510"
$ws.Range("I47").Value = "ES5566552211148855332200"
$ws.Range("I48").Value = "GB9720910936583000000000"
$ws.Range("I49").Value = "DE9301821135910205540000"
$ws.Range("I50").Value = "DE7822631245526916432102"
$ws.Range("I51").Value = "ES2120960043043075700000"
$ws.Range("I52").Value = "SM7325635478321002541225"
$ws.Range("I53").Value = "ES6832154697195423121000"
$ws.Range("I54").Value = "GR3836521452736500658485"
$ws.Range("I55").Value = "GB5520008521528775113366"
$ws.Range("I60").Value = "ES8020960043033000100000"
$ws.Range("I61").Value = "GB0836585214290025478551"
$ws.Range("I62").Value = "ES9012548523465214585214"
$ws.Range("I63").Value = "ES6931624561042546920007"
$ws.Range("I64").Value = "ES1436154231712500312566"
$ws.Range("I65").Value = "ES8244875664127231645789"
$ws.Range("I66").Value = "ES7920960031442124800000"
$ws.Range("I67").Value = "ES1633620012937852100256"
$ws.Range("I68").Value = "ES1933218885441445121022"
$ws.Range("I69").Value = "ES8462581542713690044508"
$ws.Range("I70").Value = "ES3925165151118666365100"
$ws.Range("I71").Value = "ES8020960043033000100000"
$ws.Range("I72").Value = "PT3536952365020014425254"
$ws.Range("I73").Value = "ES9565168874641561561500"
$ws.Range("I74").Value = "ES3220960583831234500000"
$ws.Range("I75").Value = "ES7221416325811510005514"
$ws.Range("I76").Value = "LU0932628484504115151115"
$ws.Range("I77").Value = "ES2220960056133231500000"
$ws.Range("I79").Value = "ES8163516541828944000984"
$ws.Range("I80").Value = "ES6223658965274585223202"
$ws.Range("I81").Value = "FI6132658012367712548745"
$ws.Range("I82").Value = "ES7223652365142254222000"
$ws.Range("I83").Value = "FR3820012541150023365233"
$ws.Range("I84").Value = "ES9232584216971684051000"
$ws.Range("I86").Value = "ES7395485212315484010000"
$ws.Range("I87").Value = "LT9321856333126985542360"
$ws.Range("I88").Value = "ES5736245978133245679001"
$ws.Range("I89").Value = "ES7631245164156597845124"
$ws.Range("I90").Value = "SM4423221158252545471411"
$ws.Range("I91").Value = "SE6832574512085411002255"
$ws.Range("I92").Value = "ES4420960043013468900000"
$ws.Range("I93").Value = "ES5631215643855060225021"
$ws.Range("I94").Value = "AT3285550564726165145610"
$ws.Range("I95").Value = "ES1665165654918886005001"
$ws.Range("I102").Value = "AT8365645150865168448896"
$ws.Range("I103").Value = "IT3526551681807651415636"
$ws.Range("I104").Value = "HU2399558741836555551120"
$ws.Range("I105").Value = "ES4352198484752100515144"
$ws.Range("I106").Value = "IE6851556584221251000254"
$ws.Range("I127").Value = "DK9032541112811220000588"
$ws.Range("I128").Value = "LT9362541122421110105611"
$ws.Range("I129").Value = "ES6855065688761051056105"
$ws.Range("I130").Value = "ES7426221011628048788896"
$ws.Range("I131").Value = "ES9712548521518742146695"
$ws.Range("I132").Value = "ES9001826530120201560000"
$ws.Range("I133").Value = "ES9021651651812511133551"
$ws.Range("I134").Value = "ES6851651487910005118185"
$ws.Range("I135").Value = "CZ9536250012804785523365"
$ws.Range("I136").Value = "AT3122515651915640081000"
